# Weekly Fruit/Vegetable price update — Mandarina (Agrícola del Norte S.A. de Arica)
#
# A new week of price records (2 rows) is inserted above the existing
# row 64, pushing all subsequent rows down by two (old row N -> new row N+2).
# The sheet's used range grows from A1:T134 to A1:T136 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 64 and 65 (existing rows 64.. shift to 66..)
$ws.Range("A64:A65").EntireRow.Insert()

# Populate both new rows with the new weekly record
for ($r = 64; $r -le 65; $r++) {
    $ws.Cells.Item($r, 1).Value  = 1
    $ws.Cells.Item($r, 2).Value  = "Agrícola del Norte S.A. de Arica"
    $ws.Cells.Item($r, 3).Value  = "Arica y Parinacota"
    $ws.Cells.Item($r, 4).Value  = 44923
    $ws.Cells.Item($r, 5).Value  = 15
    $ws.Cells.Item($r, 6).Value  = "Fruta"
    $ws.Cells.Item($r, 7).Value  = 100102
    $ws.Cells.Item($r, 8).Value  = "Cítricos"
    $ws.Cells.Item($r, 9).Value  = 100102004
    $ws.Cells.Item($r, 10).Value = "Mandarina"
    $ws.Cells.Item($r, 11).Value = "Murcott"
    $ws.Cells.Item($r, 12).Value = "Tercera"
    $ws.Cells.Item($r, 13).Value = 250
    $ws.Cells.Item($r, 14).Value = 17000
    $ws.Cells.Item($r, 15).Value = 18000
    $ws.Cells.Item($r, 16).Value = 17400
    $ws.Cells.Item($r, 17).Value = "$/caja 20 kilos"
    $ws.Cells.Item($r, 18).Value = "Región de Coquimbo"
    $ws.Cells.Item($r, 19).Value = 870
    $ws.Cells.Item($r, 20).Value = 20
}
